# Insert two new rows of weekly price data for "Brócoli" (Terminal La
# Palmera de La Serena) at row 395, pushing all the existing rows
# (previously 395-496) down to 397-498. This matches the source diff,
# which shows a new date/record pair inserted ahead of the previous
# row 395 ("Primera"/"Segunda" quality pair) while the rest of the
# table keeps its data unchanged, just shifted down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 395, shifting existing
# data (rows 395-496) down to rows 397-498.
$ws.Rows.Item(395).Resize(2).Insert()

# New row 395: "Primera" quality record for the newly added date.
$ws.Range("A395").Value = 8
$ws.Range("B395").Value = "Terminal La Palmera de La Serena"
$ws.Range("C395").Value = "Coquimbo"
$ws.Range("D395").Value = 44551
$ws.Range("E395").Value = 4
$ws.Range("F395").Value = 100112023
$ws.Range("G395").Value = "Brócoli"
$ws.Range("H395").Value = "Sin especificar"
$ws.Range("I395").Value = "Primera"
$ws.Range("J395").Value = 2200
$ws.Range("K395").Value = 600
$ws.Range("L395").Value = 700
$ws.Range("M395").Value = 650
$ws.Range("N395").Value = "$/unidad"
$ws.Range("O395").Value = "Provincia del Elquí"
$ws.Range("P395").Value = 650
$ws.Range("Q395").Value = 1
$ws.Range("R395").Value = "Hortaliza"

# New row 396: "Segunda" quality record for the newly added date.
$ws.Range("A396").Value = 8
$ws.Range("B396").Value = "Terminal La Palmera de La Serena"
$ws.Range("C396").Value = "Coquimbo"
$ws.Range("D396").Value = 44551
$ws.Range("E396").Value = 4
$ws.Range("F396").Value = 100112023
$ws.Range("G396").Value = "Brócoli"
$ws.Range("H396").Value = "Sin especificar"
$ws.Range("I396").Value = "Segunda"
$ws.Range("J396").Value = 1280
$ws.Range("K396").Value = 500
$ws.Range("L396").Value = 550
$ws.Range("M396").Value = 525
$ws.Range("N396").Value = "$/unidad"
$ws.Range("O396").Value = "Provincia del Elquí"
$ws.Range("P396").Value = 525
$ws.Range("Q396").Value = 1
$ws.Range("R396").Value = "Hortaliza"
